# "Add sample for non database" - replace the DB-backed member table sample
# with a single placeholder cell (today's date) that does not depend on a
# jdbc query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 used to hold the literal "test" label; it now holds the ${today}
# report placeholder. Escape the leading $ so PowerShell doesn't try to
# expand it as a variable/subexpression.
$ws.Range("A1").Value = "`${today}"

# The jx:each(...) comment that drove the "select * from t_member" loop
# lived on A3. Remove that comment explicitly before the cells are cleared,
# otherwise it stays anchored to the sheet.
$a3Comment = $ws.Range("A3").Comment
if ($a3Comment) {
    $a3Comment.Delete()
}

# Remove the now unused member_id / email_address header row (row 2) and
# the ${mem.member_id} / ${mem.email_address} sample row (row 3) entirely.
$ws.Range("A2:B3").EntireRow.Delete()

$wb.Save()
